# Updated capital structure database
# Applies the new column values / recalculated ratios to rows 2 and 3
# (the two Uganda "Banks (Regional)" records) on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(2, 3)

foreach ($r in $rows) {
    # Newly populated columns
    $ws.Range("D$r").Value = 0.169
    $ws.Range("E$r").Value = 0.118

    # Updated / recalculated metrics
    $ws.Range("I$r").Value = 0
    $ws.Range("J$r").Value = 0
    $ws.Range("K$r").Value = 19.9
    $ws.Range("L$r").Value = 0.2493734335839599
    $ws.Range("M$r").Value = 6.7
    $ws.Range("N$r").Value = 0.05230288836846214
    $ws.Range("O$r").Value = 0.3366834170854272
    $ws.Range("P$r").Value = 6.7
    $ws.Range("Q$r").Value = 0.05230288836846214
    $ws.Range("R$r").Value = 0.3366834170854272

    $ws.Range("U$r").Value = 94.40000000000001
    $ws.Range("V$r").Value = 0.7369242779078845
    $ws.Range("W$r").Value = 0.1417378917378917
    $ws.Range("X$r").Value = 0.08820899154041749
    $ws.Range("Y$r").Value = 0.05352890019747424
    $ws.Range("Z$r").Value = 0.5934188510875626
    $ws.Range("AA$r").Value = 0
    $ws.Range("AB$r").Value = 0.07523019489585669
    $ws.Range("AC$r").Value = -0.07523019489585669
    $ws.Range("AD$r").Value = 82.3
    $ws.Range("AE$r").Value = 0
    $ws.Range("AF$r").Value = 82.3
    $ws.Range("AG$r").Value = -12.10000000000001
    $ws.Range("AH$r").Value = 0.3911596958174905
    $ws.Range("AI$r").Value = 0.3474039679189532
    $ws.Range("AJ$r").Value = -0.1043103448275863
    $ws.Range("AK$r").Value = -0.08491228070175444

    # Columns AN and AP are no longer populated for these rows
    $ws.Range("AN$r").ClearContents()
    $ws.Range("AP$r").ClearContents()
}
